# "Corrected clock for park hours that pass midnight"
#
# The Metadata!HOURLY_PERCENT table (row 9) paired each hour-of-day with an
# arrival percentage, but the times were out of chronological order and a
# couple of late-night hours were missing altogether. This fixes the time
# sequence (10:00 ... 23:00, 00:00) and adds the missing hours, which also
# bumps TOTAL_DAILY_AGENTS and re-balances the archetype distribution to
# match, and widens a few columns so the new time values stay readable.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Metadata sheet edits
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Metadata")

# TOTAL_DAILY_AGENTS
$ws1.Range("B7").Value = 20000

# HOURLY_PERCENT: rebuild the (time, percent) pairs on row 9 in proper
# chronological order, from 10:00 through 23:00 and finally 00:00 (midnight
# the next day).
$ws1.Range("E9").Value = 15
$ws1.Range("F9").Value = 0.5
$ws1.Range("G9").Value = 12
$ws1.Range("H9").Value = 0.541666666666667
$ws1.Range("I9").Value = 13
$ws1.Range("J9").Value = 0.583333333333333
$ws1.Range("K9").Value = 8
$ws1.Range("L9").Value = 0.625
$ws1.Range("M9").Value = 10
$ws1.Range("N9").Value = 0.666666666666667
$ws1.Range("O9").Value = 13
$ws1.Range("P9").Value = 0.708333333333333
$ws1.Range("Q9").Value = 9
$ws1.Range("R9").Value = 0.75
$ws1.Range("S9").Value = 1
$ws1.Range("T9").Value = 0.791666666666667
$ws1.Range("U9").Value = 5
$ws1.Range("V9").Value = 0.833333333333333
$ws1.Range("W9").Value = 1
$ws1.Range("X9").Value = 0.875
$ws1.Range("Y9").Value = 1
$ws1.Range("Z9").Value = 0.916666666666667

# New hours added at the end of the table (22:00 was the previous end of the
# range; 23:00 and midnight were missing before).
$ws1.Range("AA9").Value = 1
$ws1.Range("AB9").Value = 0.958333333333333
$ws1.Range("AB9").NumberFormat = "hh:mm:ss\ AM/PM"
$ws1.Range("AC9").Value = 1
$ws1.Range("AD9").Value = 0
$ws1.Range("AD9").NumberFormat = "hh:mm:ss\ AM/PM"
$ws1.Range("AE9").Value = 0

# AGENT_ARCHETYPE_DISTRIBUTION: re-balance counts to match the new totals
$ws1.Range("C13").Value = 65
$ws1.Range("G13").Value = 0
$ws1.Range("I13").Value = 0

# Widen the columns used by the expanded HOURLY_PERCENT table so the extra
# time/percent pairs stay readable
$ws1.Columns.Item(2).ColumnWidth = 13
$ws1.Columns.Item(3).ColumnWidth = 3.8333333333333335
$ws1.Columns.Item(5).ColumnWidth = 2.6666666666666665
$ws1.Columns.Item(6).ColumnWidth = 10.666666666666666
$ws1.Columns.Item(7).ColumnWidth = 3.8333333333333335
$ws1.Columns.Item(8).ColumnWidth = 10.666666666666666
$ws1.Columns.Item(9).ColumnWidth = 3.8333333333333335
$ws1.Columns.Item(10).ColumnWidth = 10.666666666666666
$ws1.Columns.Item(11).ColumnWidth = 3.8333333333333335
$ws1.Columns.Item(12).ColumnWidth = 10.666666666666666
$ws1.Range("M1:AE1").EntireColumn.ColumnWidth = 3.8333333333333335

# ---------------------------------------------------------------------
# Restore the active sheet / selection state
# ---------------------------------------------------------------------
# Previously "Park" was the active tab with B4 selected; the corrected
# workbook re-opens on "Metadata" with B8 selected (just below the
# TOTAL_DAILY_AGENTS cell that was edited).
$wsActivities = $wb.Worksheets.Item("Activities")
$wsActivities.Activate()
$wsActivities.Range("B5").Select()

$ws1.Activate()
$ws1.Range("B8").Select()
